# Card22: add a new "Correction " column (N) after the existing "Event" column (M).
#   - M1 header text loses its trailing space ("Event " -> "Event")
#   - N1 gets the new header "Correction " (trailing space kept), styled like
#     the rest of the header row
#   - M2:M12 (previously blank placeholder cells) are filled with "nan",
#     matching the rest of that column's existing placeholder values
#   - N2:N12 are left blank, ready for future data entry (new column)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Card22")

# Fix the header text in M1 (drop trailing space) and add the new N1 header.
$ws.Range("M1").Value = "Event"
$ws.Range("N1").Value = "Correction "

# Give the new header cell the same look as the rest of row 1 (bold, border,
# centered) by copying the formatting from its neighbour.
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats

# Backfill the "Event" column's blank data cells with "nan", same as every
# other already-populated column in this lookup table.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"   # column M
}
